$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-24 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-25 Saturday", 2) | Out-Null
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "68÷9=7, 5", 2) | Out-Null
$d.Content.Find.Execute("45÷7=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "43÷5=8, 3", 2) | Out-Null
$d.Content.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷9=10, 7", 2) | Out-Null
$d.Content.Find.Execute("14÷5=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=7, 6", 2) | Out-Null
$d.Content.Find.Execute("89÷4=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "23÷2=11, 1", 2) | Out-Null
$d.Content.Find.Execute("58÷5=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("24÷2=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("12÷2=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷2=29, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=30, 0", 2) | Out-Null
$d.Content.Find.Execute("13÷4=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=15, 3", 2) | Out-Null
$d.Content.Find.Execute("50÷3=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2) | Out-Null
$d.Content.Find.Execute("84÷9=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "38÷7=5, 3", 2) | Out-Null
$d.Content.Find.Execute("49÷5=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "52÷8=6, 4", 2) | Out-Null
$d.Content.Find.Execute("62÷5=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷7=13, 3", 2) | Out-Null
$d.Content.Find.Execute("21÷2=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷2=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=22, 1", 2) | Out-Null
$d.Content.Find.Execute("72÷7=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=8, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2) | Out-Null
$d.Content.Find.Execute("97÷3=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=3, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷2=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "44÷8=5, 4", 2) | Out-Null
$d.Content.Find.Execute("14÷3=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null
$d.Content.Find.Execute("18÷5=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=6, 4", 2) | Out-Null
$d.Content.Find.Execute("11÷5=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "68÷8=8, 4", 2) | Out-Null
$d.Content.Find.Execute("74÷2=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷9=6, 3", 2) | Out-Null
$d.Content.Find.Execute("61÷7=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=10, 0", 2) | Out-Null

Write-Host "Replacements applied"
